# Applies the cryptos.xlsx data-refresh edit described by the commit diff.
# For D-column cells whose new text is parseable as a plain number (e.g. "0.998",
# "227.03"), a leading apostrophe forces Excel to keep/store them as literal text
# (matching the original inlineStr text cells) instead of silently converting them
# to numeric values and losing formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.226.06"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.029.74"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "'227.03"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'55.19"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'0.0787"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -4.82%  "
$ws.Range("D12").Value = "2.323.46"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'14.26"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("D14").Value = "'20.27"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "'0.744"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "'5.20"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").Value = "2.026.80"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "37.129.71"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'6.48"
$ws.Range("E19").Value = "  +7.33%  "
$ws.Range("D20").Value = "'68.89"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").Value = "'223.95"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("D26").Value = "'165.69"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "'9.22"
$ws.Range("E27").Value = "  -4.68%  "
$ws.Range("D28").Value = "'0.129"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "'18.75"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'5.56"
$ws.Range("E38").Value = "  +5.97%  "
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0215"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.470.91"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "'95.96"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'16.46"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.79"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("D45").Value = "'0.0911"
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Value = "'2.94"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "'3.66"
$ws.Range("E50").Value = "  -7.60%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.206.92"
$ws.Range("E51").Value = "  -1.24%  "
